$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 2789
$ws.Range("F4").Value = 1120
$ws.Range("F5").Value = 20332
$ws.Range("F7").Value = 2477
$ws.Range("F8").Value = 773
$ws.Range("F10").Value = 473
$ws.Range("F11").Value = 724
$ws.Range("F12").Value = 265
$ws.Range("F13").Value = 256
$ws.Range("F15").Value = 390
$ws.Range("F16").Value = 93
$ws.Range("F18").Value = 175
$ws.Range("F19").Value = 234
$ws.Range("F22").Value = 110
# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 41
$ws.Range("G2").Value = 120
$ws.Range("F6").Value = 308
$ws.Range("F12").Value = 2
$ws.Range("F15").Value = 114
$ws.Range("F23").Value = 36
# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 6063
$ws.Range("F3").Value = 672
$ws.Range("F4").Value = 637
$ws.Range("F5").Value = 1328
$ws.Range("F6").Value = 21
# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 6063
$ws.Range("F3").Value = 672
$ws.Range("F4").Value = 637
$ws.Range("F6").Value = 41
$ws.Range("G6").Value = 120
$ws.Range("F7").Value = 1328
$ws.Range("F8").Value = 2789
$ws.Range("F9").Value = 1120
$ws.Range("F10").Value = 20332
$ws.Range("F15").Value = 308
$ws.Range("F16").Value = 2477
$ws.Range("F17").Value = 773
$ws.Range("F19").Value = 21
$ws.Range("F21").Value = 473
$ws.Range("F22").Value = 724
$ws.Range("F23").Value = 265
$ws.Range("F24").Value = 256
$ws.Range("F29").Value = 390
$ws.Range("F30").Value = 93
$ws.Range("F32").Value = 2
$ws.Range("F35").Value = 175
$ws.Range("F37").Value = 234
$ws.Range("F38").Value = 114
$ws.Range("F39").Value = 114
$ws.Range("F49").Value = 36
$ws.Range("F50").Value = 110

Write-Host "Done applying updates."